# Auto-generated edit script: updates computed profit/price values per the commit diff.
# The workbook stores static numeric results (no formulas) in columns H-N of each
# job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR); this applies the refreshed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 4608.9414
$ws.Range("J32").Value2 = 5095.625
$ws.Range("L32").Value2 = 5095.625
$ws.Range("N32").Value2 = -5747.625
$ws.Range("H70").Value2 = 3480.6667
$ws.Range("J70").Value2 = 3480.6667
$ws.Range("L70").Value2 = 10442.0001
$ws.Range("N70").Value2 = -10982.0001
$ws.Range("H73").Value2 = 3480.6667
$ws.Range("J73").Value2 = 3480.6667
$ws.Range("L73").Value2 = 10442.0001
$ws.Range("N73").Value2 = -12314.0001
$ws.Range("H86").Value2 = 8873.75
$ws.Range("I86").Value2 = 9331.666999999999
$ws.Range("K86").Value2 = 9331.666999999999
$ws.Range("M86").Value2 = -8208.666999999999
$ws.Range("H89").Value2 = 8873.75
$ws.Range("I89").Value2 = 9331.666999999999
$ws.Range("K89").Value2 = 46658.335
$ws.Range("M89").Value2 = -41042.335
$ws.Range("H101").Value2 = 1263.7142
$ws.Range("J101").Value2 = 0
$ws.Range("L101").Value2 = 0
$ws.Range("N101").ClearContents()
$ws.Range("H115").Value2 = 793.4286
$ws.Range("J115").Value2 = 0
$ws.Range("L115").Value2 = 0
$ws.Range("N115").ClearContents()
$ws.Range("H127").Value2 = 1267.0667
$ws.Range("I127").Value2 = 816.3077
$ws.Range("K127").Value2 = 2448.9231
$ws.Range("M127").Value2 = 2511.0769
$ws.Range("H138").Value2 = 8271.26
$ws.Range("I138").Value2 = 6599.4
$ws.Range("J138").Value2 = 8359.253000000001
$ws.Range("K138").Value2 = 19798.2
$ws.Range("L138").Value2 = 25077.759
$ws.Range("M138").Value2 = -14658.2
$ws.Range("N138").Value2 = -35357.75900000001
$ws.Range("H141").Value2 = 1909.909
$ws.Range("J141").Value2 = 2694.3333
$ws.Range("L141").Value2 = 8082.999899999999
$ws.Range("N141").Value2 = -18442.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 19620598
$ws.Range("I32").Value2 = 23264804
$ws.Range("J32").Value2 = 32992.375
$ws.Range("K32").Value2 = 23264804
$ws.Range("L32").Value2 = 32992.375
$ws.Range("M32").Value2 = -23264517
$ws.Range("N32").Value2 = -33566.375
$ws.Range("H88").Value2 = 12544.556
$ws.Range("I88").Value2 = 21100.6
$ws.Range("J88").Value2 = 1849.5
$ws.Range("K88").Value2 = 21100.6
$ws.Range("L88").Value2 = 1849.5
$ws.Range("M88").Value2 = -20694.6
$ws.Range("N88").Value2 = -2661.5
$ws.Range("H91").Value2 = 12544.556
$ws.Range("I91").Value2 = 21100.6
$ws.Range("J91").Value2 = 1849.5
$ws.Range("K91").Value2 = 21100.6
$ws.Range("L91").Value2 = 1849.5
$ws.Range("M91").Value2 = -19696.6
$ws.Range("N91").Value2 = -4657.5
$ws.Range("H132").Value2 = 29420992
$ws.Range("I132").Value2 = 10092.931
$ws.Range("K132").Value2 = 30278.793
$ws.Range("M132").Value2 = -27748.793

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 25192.215
$ws.Range("I86").Value2 = 17130.143
$ws.Range("K86").Value2 = 17130.143
$ws.Range("M86").Value2 = -16007.143
$ws.Range("H89").Value2 = 25192.215
$ws.Range("I89").Value2 = 17130.143
$ws.Range("K89").Value2 = 85650.715
$ws.Range("M89").Value2 = -80034.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 623.2
$ws.Range("I7").Value2 = 150.83333
$ws.Range("K7").Value2 = 150.83333
$ws.Range("M7").Value2 = -37.83332999999999
$ws.Range("H31").Value2 = 26322008
$ws.Range("I31").Value2 = 4278.5713
$ws.Range("K31").Value2 = 4278.5713
$ws.Range("M31").Value2 = -3983.5713
$ws.Range("H34").Value2 = 26322008
$ws.Range("I34").Value2 = 4278.5713
$ws.Range("K34").Value2 = 4278.5713
$ws.Range("M34").Value2 = -4076.5713
$ws.Range("H132").Value2 = 83399.56
$ws.Range("I132").Value2 = 97759.42999999999
$ws.Range("J132").Value2 = 8010.25
$ws.Range("K132").Value2 = 293278.29
$ws.Range("L132").Value2 = 24030.75
$ws.Range("M132").Value2 = -290748.29
$ws.Range("N132").Value2 = -29090.75
$ws.Range("H134").Value2 = 1927.625
$ws.Range("I134").Value2 = 1900.0834
$ws.Range("J134").Value2 = 2010.25
$ws.Range("K134").Value2 = 5700.2502
$ws.Range("L134").Value2 = 6030.75
$ws.Range("M134").Value2 = -3165.2502
$ws.Range("N134").Value2 = -11100.75
$ws.Range("H141").Value2 = 324743.06
$ws.Range("J141").Value2 = 336117.66
$ws.Range("L141").Value2 = 336117.66
$ws.Range("N141").Value2 = -346477.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value2 = 211.95653
$ws.Range("I11").Value2 = 98.05
$ws.Range("K11").Value2 = 294.15
$ws.Range("M11").Value2 = -154.15
$ws.Range("H17").Value2 = 118
$ws.Range("I17").Value2 = 63.333332
$ws.Range("K17").Value2 = 189.999996
$ws.Range("M17").Value2 = -20.99999600000001
$ws.Range("H25").Value2 = 1803.1428
$ws.Range("I25").Value2 = 181.5
$ws.Range("J25").Value2 = 3965.3333
$ws.Range("K25").Value2 = 544.5
$ws.Range("L25").Value2 = 11895.9999
$ws.Range("M25").Value2 = -375.5
$ws.Range("N25").Value2 = -12233.9999
$ws.Range("H30").Value2 = 1803.1428
$ws.Range("I30").Value2 = 181.5
$ws.Range("J30").Value2 = 3965.3333
$ws.Range("K30").Value2 = 544.5
$ws.Range("L30").Value2 = 11895.9999
$ws.Range("M30").Value2 = -442.5
$ws.Range("N30").Value2 = -12099.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 6976.143
$ws.Range("J80").Value2 = 6976.143
$ws.Range("L80").Value2 = 6976.143
$ws.Range("N80").Value2 = -8972.143
$ws.Range("H83").Value2 = 6976.143
$ws.Range("J83").Value2 = 6976.143
$ws.Range("L83").Value2 = 34880.715
$ws.Range("N83").Value2 = -44864.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value2 = 1017.28
$ws.Range("J55").Value2 = 1528.5
$ws.Range("L55").Value2 = 1528.5
$ws.Range("N55").Value2 = -1874.5
$ws.Range("H100").Value2 = 4977
$ws.Range("I100").Value2 = 3802
$ws.Range("J100").Value2 = 6152
$ws.Range("K100").Value2 = 3802
$ws.Range("L100").Value2 = 6152
$ws.Range("M100").Value2 = -3261
$ws.Range("N100").Value2 = -7234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 1763.04
$ws.Range("I81").Value2 = 969.8946999999999
$ws.Range("J81").Value2 = 4274.6665
$ws.Range("K81").Value2 = 1939.7894
$ws.Range("L81").Value2 = 8549.333000000001
$ws.Range("M81").Value2 = -878.7893999999999
$ws.Range("N81").Value2 = -10671.333
$ws.Range("H84").Value2 = 1763.04
$ws.Range("I84").Value2 = 969.8946999999999
$ws.Range("J84").Value2 = 4274.6665
$ws.Range("K84").Value2 = 9698.947
$ws.Range("L84").Value2 = 42746.665
$ws.Range("M84").Value2 = -4394.947
$ws.Range("N84").Value2 = -53354.665
$ws.Range("H107").Value2 = 653.2353000000001
$ws.Range("I107").Value2 = 709.7273
$ws.Range("J107").Value2 = 549.6667
$ws.Range("K107").Value2 = 2129.1819
$ws.Range("L107").Value2 = 1649.0001
$ws.Range("M107").Value2 = -209.1819
$ws.Range("N107").Value2 = -5489.0001
$ws.Range("H113").Value2 = 850.86664
$ws.Range("I113").Value2 = 646.9
$ws.Range("J113").Value2 = 1258.8
$ws.Range("K113").Value2 = 1940.7
$ws.Range("L113").Value2 = 3776.4
$ws.Range("M113").Value2 = 229.3000000000002
$ws.Range("N113").Value2 = -8116.4
$ws.Range("H132").Value2 = 7223.3
$ws.Range("I132").Value2 = 7438.387
$ws.Range("K132").Value2 = 22315.161
$ws.Range("M132").Value2 = -19785.161
